$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Login")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Update the customer number value (C2) to a new customer number that
# doesn't result in a credit hold on orders.
$ws.Range("C2").Value = "US00000571"

# Give the updated cell a distinct 10pt black font so it stands out.
$ws.Range("C2").Font.Size = 10
$ws.Range("C2").Font.Color = 0

# Move the active cell selection to the updated cell.
$ws.Activate()
$ws.Range("C2").Select()
